# Update "想去人数" (interested-count) figures on the 展览 and 全部类型 sheets
# to reflect the newly generated output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 569   # 南宁·小蜜蜂动漫嘉年华2.0            567 -> 569
$ws1.Range("F4").Value = 355   # 南宁·漫控嘉年华09...                352 -> 355
$ws1.Range("F7").Value = 2402  # 南宁·AB动漫游戏嘉年华                2395 -> 2402
$ws1.Range("F8").Value = 413   # 横州·第二届海棠动漫游戏嘉年华         411 -> 413
$ws1.Range("F9").Value = 6233  # 南宁·第十九届（2024）良牙动漫夏季盛典  6196 -> 6233

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 569   # 南宁·小蜜蜂动漫嘉年华2.0            567 -> 569
$ws4.Range("F4").Value = 355   # 南宁·漫控嘉年华09...                352 -> 355
$ws4.Range("F9").Value = 2402  # 南宁·AB动漫游戏嘉年华                2395 -> 2402
$ws4.Range("F10").Value = 413  # 横州·第二届海棠动漫游戏嘉年华         411 -> 413
$ws4.Range("F11").Value = 6233 # 南宁·第十九届（2024）良牙动漫夏季盛典  6196 -> 6233
